# Scheduled runner update: refresh cached Leve profit calculations across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 283.63635
$ws.Range("I2").Value = 292
$ws.Range("K2").Value = 292
$ws.Range("M2").Value = -179
$ws.Range("H18").Value = 697.9
$ws.Range("I18").Value = 586.25
$ws.Range("J18").Value = 772.3333
$ws.Range("K18").Value = 586.25
$ws.Range("L18").Value = 772.3333
$ws.Range("M18").Value = -302.25
$ws.Range("N18").Value = -1340.3333
$ws.Range("H32").Value = 1375.25
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 1428.8572
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 1428.8572
$ws.Range("M32").Value = -674
$ws.Range("N32").Value = -2080.8572
$ws.Range("H38").Value = 431.41666
$ws.Range("I38").Value = 240.875
$ws.Range("J38").Value = 812.5
$ws.Range("K38").Value = 722.625
$ws.Range("L38").Value = 2437.5
$ws.Range("M38").Value = -350.625
$ws.Range("N38").Value = -3181.5
$ws.Range("H113").Value = 263832.5
$ws.Range("I113").Value = 350776.66
$ws.Range("K113").Value = 350776.66
$ws.Range("M113").Value = -347522.66
$ws.Range("H133").Value = 12579.286
$ws.Range("J133").Value = 12579.286
$ws.Range("L133").Value = 12579.286
$ws.Range("N133").Value = -22699.286
$ws.Range("H141").Value = 1897.8572
$ws.Range("I141").Value = 1742.75
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 5228.25
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -48.25
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18707.273
$ws.Range("I32").Value = 2720.6726
$ws.Range("J32").Value = 144316.28
$ws.Range("K32").Value = 2720.6726
$ws.Range("L32").Value = 144316.28
$ws.Range("M32").Value = -2433.6726
$ws.Range("N32").Value = -144890.28
$ws.Range("H61").Value = 2959.5862
$ws.Range("I61").Value = 2129.2778
$ws.Range("J61").Value = 4318.273
$ws.Range("K61").Value = 2129.2778
$ws.Range("L61").Value = 4318.273
$ws.Range("M61").Value = -1917.2778
$ws.Range("N61").Value = -4742.273
$ws.Range("H122").Value = 1924.0625
$ws.Range("I122").Value = 1546.9445
$ws.Range("J122").Value = 2408.9285
$ws.Range("K122").Value = 4640.833500000001
$ws.Range("L122").Value = 7226.7855
$ws.Range("M122").Value = -2190.833500000001
$ws.Range("N122").Value = -12126.7855
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800
$ws.Range("H136").Value = 2959.5862
$ws.Range("I136").Value = 2129.2778
$ws.Range("J136").Value = 4318.273
$ws.Range("K136").Value = 6387.8334
$ws.Range("L136").Value = 12954.819
$ws.Range("M136").Value = -3837.8334
$ws.Range("N136").Value = -18054.819
$ws.Range("H139").Value = 44457.5
$ws.Range("J139").Value = 44457.5
$ws.Range("L139").Value = 44457.5
$ws.Range("N139").Value = -54737.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 984.25
$ws.Range("I36").Value = 984.25
$ws.Range("K36").Value = 984.25
$ws.Range("M36").Value = -450.25
$ws.Range("H75").Value = 86178.25
$ws.Range("I75").Value = 6238.143
$ws.Range("K75").Value = 6238.143
$ws.Range("M75").Value = -5302.143
$ws.Range("H78").Value = 86178.25
$ws.Range("I78").Value = 6238.143
$ws.Range("K78").Value = 18714.429
$ws.Range("M78").Value = -14034.429
$ws.Range("H134").Value = 4286.7085
$ws.Range("I134").Value = 2755.75
$ws.Range("K134").Value = 8267.25
$ws.Range("M134").Value = -5732.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1485.9546
$ws.Range("I31").Value = 910.5185
$ws.Range("J31").Value = 2399.8823
$ws.Range("K31").Value = 910.5185
$ws.Range("L31").Value = 2399.8823
$ws.Range("M31").Value = -615.5185
$ws.Range("N31").Value = -2989.8823
$ws.Range("H34").Value = 1485.9546
$ws.Range("I34").Value = 910.5185
$ws.Range("J34").Value = 2399.8823
$ws.Range("K34").Value = 910.5185
$ws.Range("L34").Value = 2399.8823
$ws.Range("M34").Value = -708.5185
$ws.Range("N34").Value = -2803.8823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1897.4706
$ws.Range("I5").Value = 1521.4166
$ws.Range("J5").Value = 2800
$ws.Range("K5").Value = 4564.2498
$ws.Range("L5").Value = 8400
$ws.Range("M5").Value = -4452.2498
$ws.Range("N5").Value = -8624
$ws.Range("H80").Value = 961.63635
$ws.Range("I80").Value = 499.33334
$ws.Range("J80").Value = 1135
$ws.Range("K80").Value = 1498.00002
$ws.Range("L80").Value = 3405
$ws.Range("M80").Value = -562.0000199999999
$ws.Range("N80").Value = -5277
$ws.Range("H83").Value = 961.63635
$ws.Range("I83").Value = 499.33334
$ws.Range("J83").Value = 1135
$ws.Range("K83").Value = 4494.00006
$ws.Range("L83").Value = 10215
$ws.Range("M83").Value = 185.9999399999997
$ws.Range("N83").Value = -19575
$ws.Range("H113").Value = 12821331
$ws.Range("I113").Value = 627.0909
$ws.Range("J113").Value = 17858036
$ws.Range("K113").Value = 1881.2727
$ws.Range("L113").Value = 53574108
$ws.Range("M113").Value = 288.7273
$ws.Range("N113").Value = -53578448
$ws.Range("H132").Value = 7937785
$ws.Range("I132").Value = 750
$ws.Range("J132").Value = 9805323
$ws.Range("K132").Value = 6750
$ws.Range("L132").Value = 88247907
$ws.Range("M132").Value = -4220
$ws.Range("N132").Value = -88252967
$ws.Range("H135").Value = 1897.4706
$ws.Range("I135").Value = 1521.4166
$ws.Range("J135").Value = 2800
$ws.Range("K135").Value = 13692.7494
$ws.Range("L135").Value = 25200
$ws.Range("M135").Value = -11157.7494
$ws.Range("N135").Value = -30270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 50002.668
$ws.Range("I4").Value = 10000
$ws.Range("J4").Value = 70004
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 70004
$ws.Range("M4").Value = -9888
$ws.Range("N4").Value = -70228
$ws.Range("H122").Value = 927692.2
$ws.Range("I122").Value = 2223262
$ws.Range("K122").Value = 6669786
$ws.Range("M122").Value = -6667336
$ws.Range("H132").Value = 4498.2188
$ws.Range("I132").Value = 4573.095
$ws.Range("J132").Value = 4355.273
$ws.Range("K132").Value = 13719.285
$ws.Range("L132").Value = 13065.819
$ws.Range("M132").Value = -11189.285
$ws.Range("N132").Value = -18125.819
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 611.5714
$ws.Range("I16").Value = 611.5714
$ws.Range("K16").Value = 611.5714
$ws.Range("M16").Value = -441.5714
$ws.Range("H22").Value = 1334.2858
$ws.Range("I22").Value = 445
$ws.Range("K22").Value = 445
$ws.Range("M22").Value = -150
$ws.Range("H27").Value = 1334.2858
$ws.Range("I27").Value = 445
$ws.Range("K27").Value = 445
$ws.Range("M27").Value = -338
$ws.Range("H46").Value = 1089.1111
$ws.Range("I46").Value = 1050
$ws.Range("J46").Value = 1167.3334
$ws.Range("K46").Value = 1050
$ws.Range("L46").Value = 1167.3334
$ws.Range("M46").Value = -862
$ws.Range("N46").Value = -1543.3334
$ws.Range("H55").Value = 462.4
$ws.Range("I55").Value = 253.33333
$ws.Range("J55").Value = 776
$ws.Range("K55").Value = 253.33333
$ws.Range("L55").Value = 776
$ws.Range("M55").Value = -80.33332999999999
$ws.Range("N55").Value = -1122

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 10000
$ws.Range("J29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10580
$ws.Range("H122").Value = 35824.766
$ws.Range("I122").Value = 45348.26
$ws.Range("J122").Value = 4533.2856
$ws.Range("K122").Value = 136044.78
$ws.Range("L122").Value = 13599.8568
$ws.Range("M122").Value = -133594.78
$ws.Range("N122").Value = -18499.8568
$ws.Range("H132").Value = 7693666.5
$ws.Range("I132").Value = 10870720
$ws.Range("J132").Value = 1853.9474
$ws.Range("K132").Value = 32612160
$ws.Range("L132").Value = 5561.8422
$ws.Range("M132").Value = -32609630
$ws.Range("N132").Value = -10621.8422
